# This edit fixes formatting issues introduced when scraping floating point numbers
# (Importe column) and a handful of proveedor-name cells that used commas as separators.
# Spanish-locale formatted amounts such as "1.234,56" (dot = thousands separator,
# comma = decimal separator) are rewritten as plain "1234.56" text, and stray commas in
# a few "Razon social" values are normalized to periods, matching the source fix upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

# --- Fix proveedor names: replace "," separators with "." ---
$ws.Cells.Replace('ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN', 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN', $xlWhole)
$ws.Cells.Replace('MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO', $xlWhole)
$ws.Cells.Replace('RICCOTTI, MARIANA EDITH', 'RICCOTTI. MARIANA EDITH', $xlWhole)
$ws.Cells.Replace('SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H.', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH', $xlWhole)
$ws.Cells.Replace('NEYRA, ELISEO FELIX', 'NEYRA. ELISEO FELIX', $xlWhole)
$ws.Cells.Replace('LUGRIN, RICARDO JOSE', 'LUGRIN. RICARDO JOSE', $xlWhole)

# --- Fix Importe (amount) values: convert from Spanish-locale "1.234,56" text to plain "1234.56" text ---
# Column H ("Importe") holds these amounts as literal text. Temporarily force the column to
# Text format so Excel keeps the replaced values as text (e.g. "1900.00") instead of silently
# re-interpreting them as numbers and dropping the trailing zeros / decimal formatting.
$lastRow = $ws.UsedRange.Rows.Count
$importeRange = $ws.Range("H2:H" + $lastRow)
$importeRange.NumberFormat = "@"
$ws.Cells.Replace('1.900,00', '1900.00', $xlWhole)
$ws.Cells.Replace('1.790,00', '1790.00', $xlWhole)
$ws.Cells.Replace('158.000,00', '158000.00', $xlWhole)
$ws.Cells.Replace('332.750,00', '332750.00', $xlWhole)
$ws.Cells.Replace('495.000,00', '495000.00', $xlWhole)
$ws.Cells.Replace('5.600,00', '5600.00', $xlWhole)
$ws.Cells.Replace('717.826,20', '717826.20', $xlWhole)
$ws.Cells.Replace('3.210,00', '3210.00', $xlWhole)
$ws.Cells.Replace('9.438,00', '9438.00', $xlWhole)
$ws.Cells.Replace('90.750,00', '90750.00', $xlWhole)
$ws.Cells.Replace('7.700,00', '7700.00', $xlWhole)
$ws.Cells.Replace('700,00', '700.00', $xlWhole)
$ws.Cells.Replace('488.400,00', '488400.00', $xlWhole)
$ws.Cells.Replace('666.967,83', '666967.83', $xlWhole)
$ws.Cells.Replace('96.675,00', '96675.00', $xlWhole)
$ws.Cells.Replace('235.934,10', '235934.10', $xlWhole)
$ws.Cells.Replace('20.632,50', '20632.50', $xlWhole)
$ws.Cells.Replace('3.400,00', '3400.00', $xlWhole)
$ws.Cells.Replace('65.352,00', '65352.00', $xlWhole)
$ws.Cells.Replace('405,00', '405.00', $xlWhole)
$ws.Cells.Replace('85.393,50', '85393.50', $xlWhole)
$ws.Cells.Replace('623,00', '623.00', $xlWhole)
$ws.Cells.Replace('96.000,00', '96000.00', $xlWhole)
$ws.Cells.Replace('56.430,00', '56430.00', $xlWhole)
$ws.Cells.Replace('22.786,62', '22786.62', $xlWhole)
$ws.Cells.Replace('18.860,00', '18860.00', $xlWhole)
$ws.Cells.Replace('10.360,00', '10360.00', $xlWhole)
$ws.Cells.Replace('13.000,00', '13000.00', $xlWhole)
$ws.Cells.Replace('2.500,00', '2500.00', $xlWhole)
$ws.Cells.Replace('392,00', '392.00', $xlWhole)
$ws.Cells.Replace('25.500,00', '25500.00', $xlWhole)
$ws.Cells.Replace('40,00', '40.00', $xlWhole)
$ws.Cells.Replace('800,00', '800.00', $xlWhole)
$ws.Cells.Replace('24.500,00', '24500.00', $xlWhole)
$ws.Cells.Replace('960,40', '960.40', $xlWhole)
$ws.Cells.Replace('398.430,08', '398430.08', $xlWhole)
$ws.Cells.Replace('750,00', '750.00', $xlWhole)
$ws.Cells.Replace('102,49', '102.49', $xlWhole)
$ws.Cells.Replace('50,00', '50.00', $xlWhole)
$ws.Cells.Replace('31.600,00', '31600.00', $xlWhole)
$ws.Cells.Replace('835,60', '835.60', $xlWhole)
$ws.Cells.Replace('94,00', '94.00', $xlWhole)
$ws.Cells.Replace('17.955,51', '17955.51', $xlWhole)
$ws.Cells.Replace('72.930,00', '72930.00', $xlWhole)
$ws.Cells.Replace('687.960,00', '687960.00', $xlWhole)
$ws.Cells.Replace('37.581,52', '37581.52', $xlWhole)
$ws.Cells.Replace('4.329,54', '4329.54', $xlWhole)
$ws.Cells.Replace('213.900,00', '213900.00', $xlWhole)
$ws.Cells.Replace('399,00', '399.00', $xlWhole)
$ws.Cells.Replace('48.227,64', '48227.64', $xlWhole)
$ws.Cells.Replace('3.248,00', '3248.00', $xlWhole)
$ws.Cells.Replace('28.000,00', '28000.00', $xlWhole)
$ws.Cells.Replace('2.700,00', '2700.00', $xlWhole)
$ws.Cells.Replace('5.150,00', '5150.00', $xlWhole)
$ws.Cells.Replace('419,18', '419.18', $xlWhole)
$ws.Cells.Replace('24.649,96', '24649.96', $xlWhole)
$ws.Cells.Replace('126.500,00', '126500.00', $xlWhole)
$ws.Cells.Replace('480,20', '480.20', $xlWhole)
$ws.Cells.Replace('1.159,00', '1159.00', $xlWhole)
$ws.Cells.Replace('149,72', '149.72', $xlWhole)
$ws.Cells.Replace('13,00', '13.00', $xlWhole)
$ws.Cells.Replace('2.900,00', '2900.00', $xlWhole)
$ws.Cells.Replace('3.815,00', '3815.00', $xlWhole)
$ws.Cells.Replace('9.200,00', '9200.00', $xlWhole)
$ws.Cells.Replace('13.422,24', '13422.24', $xlWhole)
$ws.Cells.Replace('350,00', '350.00', $xlWhole)
$ws.Cells.Replace('1.312,00', '1312.00', $xlWhole)
$ws.Cells.Replace('79.685,00', '79685.00', $xlWhole)
$ws.Cells.Replace('24,90', '24.90', $xlWhole)
$ws.Cells.Replace('1.158,21', '1158.21', $xlWhole)
$ws.Cells.Replace('7.987,00', '7987.00', $xlWhole)
$ws.Cells.Replace('7.632,16', '7632.16', $xlWhole)
$ws.Cells.Replace('7.695,60', '7695.60', $xlWhole)
$ws.Cells.Replace('14.740,00', '14740.00', $xlWhole)
$ws.Cells.Replace('25.000,00', '25000.00', $xlWhole)
$ws.Cells.Replace('1.400,00', '1400.00', $xlWhole)
$ws.Cells.Replace('900,00', '900.00', $xlWhole)
$ws.Cells.Replace('57.500,00', '57500.00', $xlWhole)
$ws.Cells.Replace('3.145,00', '3145.00', $xlWhole)
$ws.Cells.Replace('726,00', '726.00', $xlWhole)
$ws.Cells.Replace('744,84', '744.84', $xlWhole)
$ws.Cells.Replace('311,40', '311.40', $xlWhole)
$ws.Cells.Replace('360,00', '360.00', $xlWhole)
$ws.Cells.Replace('7.150,00', '7150.00', $xlWhole)
$ws.Cells.Replace('10.250,00', '10250.00', $xlWhole)
$ws.Cells.Replace('190,00', '190.00', $xlWhole)
$ws.Cells.Replace('6.409,70', '6409.70', $xlWhole)
$ws.Cells.Replace('8.625,00', '8625.00', $xlWhole)
$ws.Cells.Replace('6.000,00', '6000.00', $xlWhole)
$ws.Cells.Replace('1.080,00', '1080.00', $xlWhole)
$ws.Cells.Replace('765,00', '765.00', $xlWhole)
$ws.Cells.Replace('8.200,00', '8200.00', $xlWhole)
$ws.Cells.Replace('3.600,00', '3600.00', $xlWhole)
$ws.Cells.Replace('11.230,00', '11230.00', $xlWhole)
$ws.Cells.Replace('4.200,00', '4200.00', $xlWhole)
$ws.Cells.Replace('1.925,00', '1925.00', $xlWhole)
$ws.Cells.Replace('442,89', '442.89', $xlWhole)
$ws.Cells.Replace('3.897,00', '3897.00', $xlWhole)
$ws.Cells.Replace('545,41', '545.41', $xlWhole)
$ws.Cells.Replace('43,00', '43.00', $xlWhole)
$ws.Cells.Replace('20.170,94', '20170.94', $xlWhole)
$ws.Cells.Replace('19.877,75', '19877.75', $xlWhole)
$ws.Cells.Replace('2.930,00', '2930.00', $xlWhole)
$ws.Cells.Replace('49.840,00', '49840.00', $xlWhole)
$ws.Cells.Replace('2.478,00', '2478.00', $xlWhole)
$ws.Cells.Replace('1.568,00', '1568.00', $xlWhole)
$ws.Cells.Replace('146,00', '146.00', $xlWhole)
$ws.Cells.Replace('11.419,00', '11419.00', $xlWhole)
$ws.Cells.Replace('7.498,00', '7498.00', $xlWhole)
$ws.Cells.Replace('39.824,40', '39824.40', $xlWhole)
$ws.Cells.Replace('1.694,78', '1694.78', $xlWhole)
$ws.Cells.Replace('375,00', '375.00', $xlWhole)
$ws.Cells.Replace('199,50', '199.50', $xlWhole)
$ws.Cells.Replace('10.000,00', '10000.00', $xlWhole)
$ws.Cells.Replace('630,00', '630.00', $xlWhole)
$ws.Cells.Replace('400,00', '400.00', $xlWhole)
$ws.Cells.Replace('5.030,00', '5030.00', $xlWhole)
$ws.Cells.Replace('10.300,00', '10300.00', $xlWhole)
$ws.Cells.Replace('1.263,00', '1263.00', $xlWhole)
$ws.Cells.Replace('22.000,00', '22000.00', $xlWhole)
$ws.Cells.Replace('33.900,00', '33900.00', $xlWhole)
$ws.Cells.Replace('66.020,00', '66020.00', $xlWhole)
$ws.Cells.Replace('25.533,00', '25533.00', $xlWhole)
$ws.Cells.Replace('9.469,00', '9469.00', $xlWhole)
$ws.Cells.Replace('6.540,00', '6540.00', $xlWhole)
$ws.Cells.Replace('35.428,00', '35428.00', $xlWhole)
$ws.Cells.Replace('3.500,00', '3500.00', $xlWhole)
$ws.Cells.Replace('23.000,00', '23000.00', $xlWhole)
$ws.Cells.Replace('2.300,00', '2300.00', $xlWhole)
$ws.Cells.Replace('10.800,00', '10800.00', $xlWhole)
$ws.Cells.Replace('2.168,30', '2168.30', $xlWhole)
$ws.Cells.Replace('873,90', '873.90', $xlWhole)
$ws.Cells.Replace('16.409,98', '16409.98', $xlWhole)
$ws.Cells.Replace('248,45', '248.45', $xlWhole)
$ws.Cells.Replace('1.639,00', '1639.00', $xlWhole)
$ws.Cells.Replace('3.624,35', '3624.35', $xlWhole)
$ws.Cells.Replace('4.800.000,00', '4800000.00', $xlWhole)
$ws.Cells.Replace('25.805,00', '25805.00', $xlWhole)
$ws.Cells.Replace('7.000,00', '7000.00', $xlWhole)
$ws.Cells.Replace('12.500,00', '12500.00', $xlWhole)
$ws.Cells.Replace('16.000,00', '16000.00', $xlWhole)
$ws.Cells.Replace('18.000,00', '18000.00', $xlWhole)
$ws.Cells.Replace('8.000,00', '8000.00', $xlWhole)
$ws.Cells.Replace('25.116,58', '25116.58', $xlWhole)
$ws.Cells.Replace('6.500,00', '6500.00', $xlWhole)
$ws.Cells.Replace('5.000,00', '5000.00', $xlWhole)
$ws.Cells.Replace('8.508,50', '8508.50', $xlWhole)
$ws.Cells.Replace('9.000,00', '9000.00', $xlWhole)
$ws.Cells.Replace('14.000,00', '14000.00', $xlWhole)
$ws.Cells.Replace('17.500,00', '17500.00', $xlWhole)
$ws.Cells.Replace('24.000,00', '24000.00', $xlWhole)
$ws.Cells.Replace('24.860,00', '24860.00', $xlWhole)
$ws.Cells.Replace('15.000,00', '15000.00', $xlWhole)
$ws.Cells.Replace('9.910,00', '9910.00', $xlWhole)
$ws.Cells.Replace('39.740,00', '39740.00', $xlWhole)
$ws.Cells.Replace('4.000,00', '4000.00', $xlWhole)
$ws.Cells.Replace('9.130,00', '9130.00', $xlWhole)
$ws.Cells.Replace('510,00', '510.00', $xlWhole)
$ws.Cells.Replace('162.500,00', '162500.00', $xlWhole)
$ws.Cells.Replace('2.100,00', '2100.00', $xlWhole)
$ws.Cells.Replace('234.550,00', '234550.00', $xlWhole)
$ws.Cells.Replace('9.900,00', '9900.00', $xlWhole)
$ws.Cells.Replace('2.204,00', '2204.00', $xlWhole)
$ws.Cells.Replace('4.848,00', '4848.00', $xlWhole)
$ws.Cells.Replace('30.682,10', '30682.10', $xlWhole)
$ws.Cells.Replace('18.365,00', '18365.00', $xlWhole)
$ws.Cells.Replace('139,00', '139.00', $xlWhole)
$ws.Cells.Replace('16.720,00', '16720.00', $xlWhole)
$ws.Cells.Replace('3.447,00', '3447.00', $xlWhole)
$ws.Cells.Replace('957,56', '957.56', $xlWhole)
$ws.Cells.Replace('8.835,00', '8835.00', $xlWhole)
$ws.Cells.Replace('4.100,50', '4100.50', $xlWhole)
$ws.Cells.Replace('5.616,06', '5616.06', $xlWhole)
$ws.Cells.Replace('5.400,00', '5400.00', $xlWhole)
$ws.Cells.Replace('1.010,20', '1010.20', $xlWhole)
$ws.Cells.Replace('1.704,00', '1704.00', $xlWhole)
$ws.Cells.Replace('220,00', '220.00', $xlWhole)
$ws.Cells.Replace('36.834,68', '36834.68', $xlWhole)
$ws.Cells.Replace('1.802,12', '1802.12', $xlWhole)
$ws.Cells.Replace('7.800,00', '7800.00', $xlWhole)
$ws.Cells.Replace('28.342,00', '28342.00', $xlWhole)
$ws.Cells.Replace('6.201,04', '6201.04', $xlWhole)
$ws.Cells.Replace('7.900,00', '7900.00', $xlWhole)
$ws.Cells.Replace('59,74', '59.74', $xlWhole)
$ws.Cells.Replace('2.850,00', '2850.00', $xlWhole)
$ws.Cells.Replace('2.000,00', '2000.00', $xlWhole)
$ws.Cells.Replace('156,74', '156.74', $xlWhole)
$ws.Cells.Replace('7.081,64', '7081.64', $xlWhole)
$ws.Cells.Replace('82.800,00', '82800.00', $xlWhole)
$ws.Cells.Replace('30.000,00', '30000.00', $xlWhole)
$ws.Cells.Replace('77.700,00', '77700.00', $xlWhole)
$ws.Cells.Replace('60.000,00', '60000.00', $xlWhole)
$ws.Cells.Replace('64.500,00', '64500.00', $xlWhole)
$ws.Cells.Replace('15.802,70', '15802.70', $xlWhole)
$ws.Cells.Replace('4.413.721,18', '4413721.18', $xlWhole)
$ws.Cells.Replace('500,00', '500.00', $xlWhole)
$ws.Cells.Replace('6.450,00', '6450.00', $xlWhole)
$ws.Cells.Replace('14.400,00', '14400.00', $xlWhole)
$ws.Cells.Replace('150.500,00', '150500.00', $xlWhole)
$ws.Cells.Replace('159.500,00', '159500.00', $xlWhole)
$ws.Cells.Replace('152.750,00', '152750.00', $xlWhole)
$ws.Cells.Replace('269.500,00', '269500.00', $xlWhole)
$ws.Cells.Replace('337.500,00', '337500.00', $xlWhole)
$ws.Cells.Replace('394.500,00', '394500.00', $xlWhole)
$ws.Cells.Replace('167.500,00', '167500.00', $xlWhole)
$ws.Cells.Replace('388.500,00', '388500.00', $xlWhole)
$ws.Cells.Replace('275.687,50', '275687.50', $xlWhole)
$ws.Cells.Replace('278.000,00', '278000.00', $xlWhole)
$ws.Cells.Replace('155.250,00', '155250.00', $xlWhole)
$ws.Cells.Replace('264.823,00', '264823.00', $xlWhole)
$ws.Cells.Replace('1.036.654,62', '1036654.62', $xlWhole)
$ws.Cells.Replace('15.900,00', '15900.00', $xlWhole)
$ws.Cells.Replace('35.000,00', '35000.00', $xlWhole)
$ws.Cells.Replace('256.224,48', '256224.48', $xlWhole)
$ws.Cells.Replace('155.800,00', '155800.00', $xlWhole)
$ws.Cells.Replace('9.300,00', '9300.00', $xlWhole)
$ws.Cells.Replace('20.000,00', '20000.00', $xlWhole)
$ws.Cells.Replace('77.850,00', '77850.00', $xlWhole)
$ws.Cells.Replace('86.230,00', '86230.00', $xlWhole)
$ws.Cells.Replace('264.000,00', '264000.00', $xlWhole)
$ws.Cells.Replace('70.000,00', '70000.00', $xlWhole)
$ws.Cells.Replace('1.551.000,00', '1551000.00', $xlWhole)
$ws.Cells.Replace('9.100,00', '9100.00', $xlWhole)
$ws.Cells.Replace('5.800,00', '5800.00', $xlWhole)
$ws.Cells.Replace('4.409,31', '4409.31', $xlWhole)
$ws.Cells.Replace('8.668,44', '8668.44', $xlWhole)
$ws.Cells.Replace('4.100,00', '4100.00', $xlWhole)
$ws.Cells.Replace('12.155,00', '12155.00', $xlWhole)
$ws.Cells.Replace('170.856,00', '170856.00', $xlWhole)
$ws.Cells.Replace('2.880,00', '2880.00', $xlWhole)
$ws.Cells.Replace('6.600,00', '6600.00', $xlWhole)
$ws.Cells.Replace('960,00', '960.00', $xlWhole)
$ws.Cells.Replace('1.800,00', '1800.00', $xlWhole)

# Restore the default ("Normal") style so the cell formatting matches the original workbook.
$importeRange.Style = "Normal"
